# Update workbook with the latest monthly data point (01-08-2021)
# and revised figures for 01-05-2021, 01-06-2021 and 01-07-2021.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to existing rows -----------------------------------

# Row 174 -> 01-05-2021
$ws.Cells.Item(174, 2).Value  = 6681
$ws.Cells.Item(174, 3).Value  = 6681
$ws.Cells.Item(174, 5).Value  = 2167
$ws.Cells.Item(174, 8).Value  = 2603
$ws.Cells.Item(174, 9).Value  = 566
$ws.Cells.Item(174, 10).Value = 341

# Row 175 -> 01-06-2021
$ws.Cells.Item(175, 8).Value = 2727
$ws.Cells.Item(175, 9).Value = 731

# Row 176 -> 01-07-2021
$ws.Cells.Item(176, 2).Value  = 6246
$ws.Cells.Item(176, 3).Value  = 6246
$ws.Cells.Item(176, 4).Value  = 4159
$ws.Cells.Item(176, 5).Value  = 2087
$ws.Cells.Item(176, 7).Value  = 1622
$ws.Cells.Item(176, 8).Value  = 2701
$ws.Cells.Item(176, 9).Value  = 626
$ws.Cells.Item(176, 10).Value = 298
$ws.Cells.Item(176, 11).Value = 112

# --- New row for 01-08-2021 ------------------------------------------

$ws.Cells.Item(177, 1).NumberFormat = "@"
$ws.Cells.Item(177, 1).Value  = "01-08-2021"
$ws.Cells.Item(177, 1).ClearFormats()
$ws.Cells.Item(177, 2).Value  = 6918
$ws.Cells.Item(177, 3).Value  = 6918
$ws.Cells.Item(177, 4).Value  = 4680
$ws.Cells.Item(177, 5).Value  = 2239
$ws.Cells.Item(177, 6).Value  = 842
$ws.Cells.Item(177, 7).Value  = 1813
$ws.Cells.Item(177, 8).Value  = 2682
$ws.Cells.Item(177, 9).Value  = 594
$ws.Cells.Item(177, 10).Value = 374
$ws.Cells.Item(177, 11).Value = 142
$ws.Cells.Item(177, 12).Value = 471
$ws.Cells.Item(177, 13).Value = 0
$ws.Cells.Item(177, 14).Value = 0
$ws.Cells.Item(177, 15).Value = 0
